$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before F (old F -> G, shifting addresses/district data right)
$ws.Columns("F:F").Insert()

# Populate the new Address column (F) for rows where the address is known
$ws.Range("F2").Value = "Address"
$ws.Range("F3").Value = "S K R HS GL HalliSondur"
$ws.Range("F5").Value = "S M H S GanadahunaseSira"
$ws.Range("F11").Value = "S G R H S Kote Hosadurga"
$ws.Range("F12").Value = "Sri N Ghattappa nayaka Rural HS ChickkobanahalliMolakalmuru"
$ws.Range("F13").Value = "Adarsha Vidyalaya hagaribommanahalli"
$ws.Range("F15").Value = "M G V P H S CB Kere Hosadurga"
$ws.Range("F17").Value = "Shanthala High School Chikkanahally Sira"
$ws.Range("F19").Value = "GHS MankiHonnavar"
$ws.Range("F20").Value = "S B R H S BG KereMolakalmuru"
$ws.Range("F24").Value = "Govt. Adarshavidyalaya Siruguppa"
$ws.Range("F27").Value = "Sree Basaveswar Rural HS Mathodu Hosadurga"
$ws.Range("F28").Value = "GHS Kurubarahalli"
$ws.Range("F31").Value = "Jeevanidhi High School Bukkapattana Sira"
$ws.Range("F32").Value = "GHS RavihalSiruguppa"
$ws.Range("F35").Value = "S K.D.D.V HS"
$ws.Range("F36").Value = "S S R C H S Vijapura"
$ws.Range("F39").Value = "Govt. Adarsha Vidyalaya Hagaribommanahalli"
$ws.Range("F40").Value = "GHS MuddatanurSiruguppa"
